$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.613.96"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").Value = "3.064.58"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.28"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.69"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.051.67"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +3.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  +9.31%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.81"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "3.564.32"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "3.058.70"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").Value = "61.560.66"
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.56"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.727"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.68"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.22"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.99"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  +6.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.44"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("E33").Value = "  +8.50%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.03"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0790"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.04"
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.16"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.00"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.80"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "420.68"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0366"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").Value = "2.762.57"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  +6.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.83"
$ws.Range("E46").Value = "  +8.53%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.12"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.74"
$ws.Range("E51").Value = "  +0.39%  "
